{"js": "// Replace the outgoing HSE manager's signature block with the new\n// technician's name/title in the NR-12 authorization signature table.\n//\n// Two paragraphs are touched (both in the left column of the second\n// signature-table row):\n//   1) \"Bruna Petroni Ces\u00e1rio\"   -> \"LEONARDO SILVERIO FERREIRA\"\n//      (also drops the \"TableParagraph\" paragraph style and inlines its\n//      font/language properties onto the paragraph mark + run instead)\n//   2) \"Gerente de HSE Brasil\"   -> \"T\u00e9cnico(a) de Seguran\u00e7a do Trabalho\"\n//      (also adds the eastAsia font + language tags to the run)\n\nasync function replaceParagraphOoxml(context, searchText, ooxmlBody) {\n  const results = context.document.body.search(searchText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n\n  const para = results.items[0].paragraphs.getFirst();\n  const range = para.getRange(\"Whole\");\n\n  const ooxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + ooxmlBody + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>';\n\n  range.insertOoxml(ooxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Name paragraph: drop pStyle=\"TableParagraph\", inline the Arial\n//    rFonts + pt-PT lang on both the paragraph mark and the run, and\n//    swap the name text. Bold stays on both rPr's.\nawait replaceParagraphOoxml(\n  context,\n  \"Bruna Petroni Ces\u00e1rio\",\n  '<w:p>' +\n    '<w:pPr>' +\n      '<w:jc w:val=\"center\"/>' +\n      '<w:rPr>' +\n        '<w:rFonts w:ascii=\"Arial\" w:eastAsia=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n        '<w:b/><w:bCs/>' +\n        '<w:lang w:eastAsia=\"pt-PT\" w:bidi=\"pt-PT\"/>' +\n      '</w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:ascii=\"Arial\" w:eastAsia=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n        '<w:b/><w:bCs/>' +\n        '<w:lang w:eastAsia=\"pt-PT\" w:bidi=\"pt-PT\"/>' +\n      '</w:rPr>' +\n      '<w:t>LEONARDO SILVERIO FERREIRA</w:t>' +\n    '</w:r>' +\n  '</w:p>'\n);\n\n// 2) Title paragraph: add eastAsia Arial + lang to the run, swap text.\nawait replaceParagraphOoxml(\n  context,\n  \"Gerente de HSE Brasil\",\n  '<w:p>' +\n    '<w:pPr>' +\n      '<w:jc w:val=\"center\"/>' +\n      '<w:rPr>' +\n        '<w:rFonts w:ascii=\"Arial\" w:eastAsia=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n        '<w:lang w:val=\"pt-PT\" w:eastAsia=\"pt-PT\" w:bidi=\"pt-PT\"/>' +\n      '</w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n      '<w:rPr>' +\n        '<w:rFonts w:ascii=\"Arial\" w:eastAsia=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n        '<w:lang w:eastAsia=\"pt-PT\" w:bidi=\"pt-PT\"/>' +\n      '</w:rPr>' +\n      '<w:t>T\u00e9cnico(a) de Seguran\u00e7a do Trabalho</w:t>' +\n    '</w:r>' +\n  '</w:p>'\n);\n", "ps1": "# Replace the outgoing HSE manager's signature block with the new\n# technician's name/title in the NR-12 authorization signature table.\n#\n# Two paragraphs are touched (both in the left column of the second\n# signature-table row):\n#   1) \"Bruna Petroni Ces\u00e1rio\"   -> \"LEONARDO SILVERIO FERREIRA\"\n#      (also drops the \"TableParagraph\" paragraph style and inlines its\n#      font/language properties onto the paragraph mark + run instead)\n#   2) \"Gerente de HSE Brasil\"   -> \"T\u00e9cnico(a) de Seguran\u00e7a do Trabalho\"\n#      (also adds the eastAsia font + language tags to the run)\n\n$d = $word.ActiveDocument\n\nfunction Replace-ParagraphOoxml($doc, [string]$searchText, [string]$ooxmlBody) {\n    $rng = $doc.Content\n    $find = $rng.Find\n    $find.Text = $searchText\n    $find.MatchCase = $true\n    $find.Execute() | Out-Null\n\n    $para = $rng.Paragraphs(1)\n    $pr = $para.Range\n\n    $ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $ooxmlBody + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n\n    $pr.InsertXML($ooxml) | Out-Null\n}\n\n# 1) Name paragraph: drop pStyle=\"TableParagraph\", inline the Arial\n#    rFonts + pt-PT lang on both the paragraph mark and the run, and\n#    swap the name text. Bold stays on both rPr's.\n$nameBody = '<w:p>' +\n    '<w:pPr>' +\n        '<w:jc w:val=\"center\"/>' +\n        '<w:rPr>' +\n            '<w:rFonts w:ascii=\"Arial\" w:eastAsia=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n            '<w:b/><w:bCs/>' +\n            '<w:lang w:eastAsia=\"pt-PT\" w:bidi=\"pt-PT\"/>' +\n        '</w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n        '<w:rPr>' +\n            '<w:rFonts w:ascii=\"Arial\" w:eastAsia=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n            '<w:b/><w:bCs/>' +\n            '<w:lang w:eastAsia=\"pt-PT\" w:bidi=\"pt-PT\"/>' +\n        '</w:rPr>' +\n        '<w:t>LEONARDO SILVERIO FERREIRA</w:t>' +\n    '</w:r>' +\n'</w:p>'\nReplace-ParagraphOoxml $d \"Bruna Petroni Ces\u00e1rio\" $nameBody\n\n# 2) Title paragraph: add eastAsia Arial + lang to the run, swap text.\n$titleBody = '<w:p>' +\n    '<w:pPr>' +\n        '<w:jc w:val=\"center\"/>' +\n        '<w:rPr>' +\n            '<w:rFonts w:ascii=\"Arial\" w:eastAsia=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n            '<w:lang w:val=\"pt-PT\" w:eastAsia=\"pt-PT\" w:bidi=\"pt-PT\"/>' +\n        '</w:rPr>' +\n    '</w:pPr>' +\n    '<w:r>' +\n        '<w:rPr>' +\n            '<w:rFonts w:ascii=\"Arial\" w:eastAsia=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/>' +\n            '<w:lang w:eastAsia=\"pt-PT\" w:bidi=\"pt-PT\"/>' +\n        '</w:rPr>' +\n        '<w:t>T\u00e9cnico(a) de Seguran\u00e7a do Trabalho</w:t>' +\n    '</w:r>' +\n'</w:p>'\nReplace-ParagraphOoxml $d \"Gerente de HSE Brasil\" $titleBody\n"}
